$d = $word.ActiveDocument

# --- 1. Locate the "Once you have finished..." paragraph and append the new   ---
#        sentence to it (as a trailing run).
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Once you have finished the game you will be taken to the following end screen where you will see your score.",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    # Find.Execute narrows $searchRange to the matched text; collapse to its end
    # and insert the additional sentence as new run content right after it.
    $searchRange.Collapse(0)
    [void]$searchRange.InsertAfter(" And if it has been a high score. You can then squeeze again to go back to the main menu and if you have set a high score you can see it displayed on the main menu.")
}

# --- 2. Insert a new paragraph "Thank you for watching..." right after the    ---
#        blank paragraph that follows the paragraph above, leaving the other
#        blank paragraphs untouched.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Once you have finished the game you will be taken to the following end screen*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $blankPara = $d.Paragraphs.Item($targetIndex + 1)
    $insertionPoint = $blankPara.Range
    $insertionPoint.Collapse(0)
    [void]$insertionPoint.InsertParagraphAfter()

    $thankYouPara = $d.Paragraphs.Item($targetIndex + 2)
    $thankYouPara.Range.Text = "Thank you for watching if you have any questions please let me know"
}
